$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 75 (Leve Item ID context G75)
$ws.Range("H75").Value = 29950
$ws.Range("J75").Value = 29950
$ws.Range("L75").Value = 29950
$ws.Range("N75").Value = -31822

# Row 78 (Leve Item ID context G78)
$ws.Range("H78").Value = 29950
$ws.Range("J78").Value = 29950
$ws.Range("L78").Value = 89850
$ws.Range("N78").Value = -99210

# Row 107 (Leve Item ID context G107)
$ws.Range("H107").Value = 200.6
$ws.Range("I107").Value = 200.6
$ws.Range("K107").Value = 200.6
$ws.Range("M107").Value = 1719.4

# Row 111 (Leve Item ID context G111)
$ws.Range("H111").Value = 2003
$ws.Range("I111").Value = 1362.875
$ws.Range("J111").Value = 3466.1428
$ws.Range("K111").Value = 4088.625
$ws.Range("L111").Value = 10398.4284
$ws.Range("M111").Value = -1021.625
$ws.Range("N111").Value = -16532.4284

# Row 112 (Leve Item ID context G112)
$ws.Range("H112").Value = 1392.8572
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1392.8572
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4178.571599999999
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6394.571599999999

# Row 113 (Leve Item ID context G113)
$ws.Range("H113").Value = 4252.7646
$ws.Range("I113").Value = 3699.875
$ws.Range("J113").Value = 4744.222
$ws.Range("K113").Value = 3699.875
$ws.Range("L113").Value = 4744.222
$ws.Range("M113").Value = -445.875
$ws.Range("N113").Value = -11252.222

# Row 115 (Leve Item ID context G115)
$ws.Range("H115").Value = 1721.4706
$ws.Range("J115").Value = 3480
$ws.Range("L115").Value = 10440
$ws.Range("N115").Value = -13574

# Row 116 (Leve Item ID context G116)
$ws.Range("H116").Value = 125002750
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 250002500
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 250002500
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -250009384

# Row 118 (Leve Item ID context G118)
$ws.Range("H118").Value = 762.25
$ws.Range("I118").Value = 442.57144
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 1327.71432
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 329.28568
$ws.Range("N118").Value = -12314

# Row 121 (Leve Item ID context G121)
$ws.Range("H121").Value = 1105
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1105
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3315
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -6809

# Row 129 (Leve Item ID context G129)
$ws.Range("H129").Value = 1010.4375
$ws.Range("J129").Value = 1188.3334
$ws.Range("L129").Value = 3565.0002
$ws.Range("N129").Value = -13565.0002

# Row 137 (Leve Item ID context G137)
$ws.Range("H137").Value = 1986729.1
$ws.Range("I137").Value = 2512.1333
$ws.Range("K137").Value = 7536.3999
$ws.Range("M137").Value = -4986.3999

# Row 138 (Leve Item ID context G138)
$ws.Range("H138").Value = 7696187
$ws.Range("I138").Value = 2590.6
$ws.Range("J138").Value = 12504685
$ws.Range("K138").Value = 7771.799999999999
$ws.Range("L138").Value = 37514055
$ws.Range("M138").Value = -2631.799999999999
$ws.Range("N138").Value = -37524335


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID context G2)
$ws.Range("H2").Value = 1524.5714
$ws.Range("I2").Value = 1140
$ws.Range("J2").Value = 2037.3334
$ws.Range("K2").Value = 1140
$ws.Range("L2").Value = 2037.3334
$ws.Range("M2").Value = -1027
$ws.Range("N2").Value = -2263.3334

# Row 110 (Leve Item ID context G110)
$ws.Range("H110").Value = 1043.0416
$ws.Range("I110").Value = 1047.8636
$ws.Range("J110").Value = 990
$ws.Range("K110").Value = 1047.8636
$ws.Range("L110").Value = 990
$ws.Range("M110").Value = 997.1364000000001
$ws.Range("N110").Value = -5080

# Row 116 (Leve Item ID context G116)
$ws.Range("H116").Value = 1524.5714
$ws.Range("I116").Value = 1140
$ws.Range("J116").Value = 2037.3334
$ws.Range("K116").Value = 1140
$ws.Range("L116").Value = 2037.3334
$ws.Range("M116").Value = 1154
$ws.Range("N116").Value = -6625.3334


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID context G3)
$ws.Range("H3").Value = 1524.5714
$ws.Range("I3").Value = 1140
$ws.Range("J3").Value = 2037.3334
$ws.Range("K3").Value = 1140
$ws.Range("L3").Value = 2037.3334
$ws.Range("M3").Value = -1026
$ws.Range("N3").Value = -2265.3334

# Row 20 (Leve Item ID context G20)
$ws.Range("H20").Value = 1272.909
$ws.Range("I20").Value = 786.2727
$ws.Range("J20").Value = 1516.2273
$ws.Range("K20").Value = 786.2727
$ws.Range("L20").Value = 1516.2273
$ws.Range("M20").Value = -539.2727
$ws.Range("N20").Value = -2010.2273

# Row 37 (Leve Item ID context G37)
$ws.Range("H37").Value = 10978.333
$ws.Range("I37").Value = 2495
$ws.Range("J37").Value = 17765
$ws.Range("K37").Value = 2495
$ws.Range("L37").Value = 17765
$ws.Range("M37").Value = -2358
$ws.Range("N37").Value = -18039

# Row 86 (Leve Item ID context G86)
$ws.Range("H86").Value = 1810.1515
$ws.Range("I86").Value = 1785.5454
$ws.Range("J86").Value = 1859.3636
$ws.Range("K86").Value = 1785.5454
$ws.Range("L86").Value = 1859.3636
$ws.Range("M86").Value = -662.5454
$ws.Range("N86").Value = -4105.3636

# Row 88 (Leve Item ID context G88)
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 89 (Leve Item ID context G89)
$ws.Range("H89").Value = 1810.1515
$ws.Range("I89").Value = 1785.5454
$ws.Range("J89").Value = 1859.3636
$ws.Range("K89").Value = 8927.726999999999
$ws.Range("L89").Value = 9296.817999999999
$ws.Range("M89").Value = -3311.726999999999
$ws.Range("N89").Value = -20528.818

# Row 91 (Leve Item ID context G91)
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# Row 107 (Leve Item ID context G107)
$ws.Range("H107").Value = 2731.125
$ws.Range("I107").Value = 2391.3684
$ws.Range("J107").Value = 4022.2
$ws.Range("K107").Value = 2391.3684
$ws.Range("L107").Value = 4022.2
$ws.Range("M107").Value = -471.3683999999998
$ws.Range("N107").Value = -7862.2

# Row 134 (Leve Item ID context G134)
$ws.Range("H134").Value = 6914.875
$ws.Range("I134").Value = 7496.647
$ws.Range("J134").Value = 5502
$ws.Range("K134").Value = 22489.941
$ws.Range("L134").Value = 16506
$ws.Range("M134").Value = -19954.941
$ws.Range("N134").Value = -21576


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID context G31)
$ws.Range("H31").Value = 4132.4375
$ws.Range("I31").Value = 2378.0908
$ws.Range("J31").Value = 7992
$ws.Range("K31").Value = 2378.0908
$ws.Range("L31").Value = 7992
$ws.Range("M31").Value = -2083.0908
$ws.Range("N31").Value = -8582

# Row 34 (Leve Item ID context G34)
$ws.Range("H34").Value = 4132.4375
$ws.Range("I34").Value = 2378.0908
$ws.Range("J34").Value = 7992
$ws.Range("K34").Value = 2378.0908
$ws.Range("L34").Value = 7992
$ws.Range("M34").Value = -2176.0908
$ws.Range("N34").Value = -8396

# Row 122 (Leve Item ID context G122)
$ws.Range("H122").Value = 1152.6666
$ws.Range("I122").Value = 959
$ws.Range("J122").Value = 1540
$ws.Range("K122").Value = 2877
$ws.Range("L122").Value = 4620
$ws.Range("M122").Value = -427
$ws.Range("N122").Value = -9520


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID context G5)
$ws.Range("H5").Value = 732.375
$ws.Range("I5").Value = 402.1111
$ws.Range("K5").Value = 1206.3333
$ws.Range("M5").Value = -1094.3333

# Row 120 (Leve Item ID context G120)
$ws.Range("H120").Value = 8792
$ws.Range("I120").Value = 4000
$ws.Range("J120").Value = 11986.667
$ws.Range("K120").Value = 12000
$ws.Range("L120").Value = 35960.001
$ws.Range("M120").Value = -7162
$ws.Range("N120").Value = -45636.001

# Row 122 (Leve Item ID context G122)
$ws.Range("H122").Value = 1470.7333
$ws.Range("I122").Value = 994.25
$ws.Range("J122").Value = 1644
$ws.Range("K122").Value = 8948.25
$ws.Range("L122").Value = 14796
$ws.Range("M122").Value = -6498.25
$ws.Range("N122").Value = -19696

# Row 131 (Leve Item ID context G131)
$ws.Range("H131").Value = 887.4
$ws.Range("I131").Value = 586.6667
$ws.Range("J131").Value = 896.70105
$ws.Range("K131").Value = 1760.0001
$ws.Range("L131").Value = 2690.10315
$ws.Range("M131").Value = 3279.9999
$ws.Range("N131").Value = -12770.10315

# Row 135 (Leve Item ID context G135)
$ws.Range("H135").Value = 732.375
$ws.Range("I135").Value = 402.1111
$ws.Range("K135").Value = 3618.9999
$ws.Range("M135").Value = -1083.9999


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID context G70)
$ws.Range("H70").Value = 4806.354
$ws.Range("I70").Value = 4672.9644
$ws.Range("J70").Value = 5636.3335
$ws.Range("K70").Value = 4672.9644
$ws.Range("L70").Value = 5636.3335
$ws.Range("M70").Value = -4402.9644
$ws.Range("N70").Value = -6176.3335

# Row 73 (Leve Item ID context G73)
$ws.Range("H73").Value = 4806.354
$ws.Range("I73").Value = 4672.9644
$ws.Range("J73").Value = 5636.3335
$ws.Range("K73").Value = 4672.9644
$ws.Range("L73").Value = 5636.3335
$ws.Range("M73").Value = -3736.9644
$ws.Range("N73").Value = -7508.3335

# Row 113 (Leve Item ID context G113)
$ws.Range("H113").Value = 1622.5862
$ws.Range("I113").Value = 1294.8462
$ws.Range("J113").Value = 1888.875
$ws.Range("K113").Value = 1294.8462
$ws.Range("L113").Value = 1888.875
$ws.Range("M113").Value = 875.1538
$ws.Range("N113").Value = -6228.875


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID context G22)
$ws.Range("H22").Value = 911.6667
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 911.6667
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 911.6667
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1501.6667

# Row 27 (Leve Item ID context G27)
$ws.Range("H27").Value = 911.6667
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 911.6667
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 911.6667
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1125.6667

# Row 61 (Leve Item ID context G61)
$ws.Range("H61").Value = 1999.8422
$ws.Range("I61").Value = 1682.6666
$ws.Range("J61").Value = 2543.5715
$ws.Range("K61").Value = 1682.6666
$ws.Range("L61").Value = 2543.5715
$ws.Range("M61").Value = -1480.6666
$ws.Range("N61").Value = -2947.5715

# Row 113 (Leve Item ID context G113)
$ws.Range("H113").Value = 1999.8422
$ws.Range("I113").Value = 1682.6666
$ws.Range("J113").Value = 2543.5715
$ws.Range("K113").Value = 1682.6666
$ws.Range("L113").Value = 2543.5715
$ws.Range("M113").Value = 487.3334
$ws.Range("N113").Value = -6883.5715

